$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.502.08"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.838.87"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.22"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5380"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2967"
$ws.Range("E8").Value = "  -8.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06957"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.46"
$ws.Range("E10").Value = "  -7.90%  "
$ws.Range("D11").Value = "1.856.38"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7289"
$ws.Range("E12").Value = "  -6.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07223"
$ws.Range("E13").Value = "  -6.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.18"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007894"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "26.515.38"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "2.080.54"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.587"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.999"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.206"
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.30"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.167"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.695"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.99"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.21"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.243"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08880"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.038"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04845"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.916"
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7244"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.089"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.299"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01712"
$ws.Range("E39").Value = "  -4.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4696"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9041"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.15"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.405"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1248"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.000"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4064"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.78"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8919"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05758"
$ws.Range("E51").Value = "  -2.08%  "
